# Adds "order" and "color" fields (two new columns) to the command deck,
# inserted immediately before the existing "example"/"description" columns
# (old AB/AC), shifting those to AD/AE.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two new blank columns at AB:AC ------------------------------
$ws.Range("AB1:AC1").EntireColumn.Insert()

# --- 2. Header row 1 (merged-cell spacer row): style only ------------------
$ws.Range("L2").Copy()
$ws.Range("AB1:AC1").PasteSpecial(-4122)   # xlPasteFormats

# --- 3. Header row 2: "order" / "color" labels -----------------------------
$ws.Range("A2").Copy()
$ws.Range("AB2:AC2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(2, 28).Value = "order"
$ws.Cells.Item(2, 29).Value = "color"

# --- 4. Spacer row 33: blank, same look as column AA -----------------------
$ws.Range("AA33").Copy()
$ws.Range("AB33:AC33").PasteSpecial(-4122) # xlPasteFormats

# --- 5. Data rows: numeric order + color hex string ------------------------
$ws.Range("A3").Copy()
$ws.Range("AB3:AC32").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("AB34:AC85").PasteSpecial(-4122) # xlPasteFormats

$order = @{
  3=71; 4=72; 5=73; 6=44; 7=45; 8=46; 9=47; 10=48; 11=49; 12=50;
  13=51; 14=52; 15=53; 16=54; 17=55; 18=56; 19=57; 20=58; 21=59; 22=60;
  23=61; 24=62; 25=63; 26=64; 27=65; 28=66; 29=67; 30=68; 31=69; 32=70;
  34=74; 35=75; 36=7; 37=8; 38=9; 39=82; 40=10; 41=2; 42=1; 43=3;
  44=5; 45=6; 46=4; 47=11; 48=12; 49=13; 50=14; 51=15; 52=16; 53=17;
  54=18; 55=19; 56=20; 57=21; 58=22; 59=23; 60=24; 61=25; 62=26; 63=27;
  64=28; 65=29; 66=30; 67=31; 68=32; 69=33; 70=34; 71=35; 72=36; 73=37;
  74=38; 75=39; 76=40; 77=41; 78=42; 79=43; 80=76; 81=77; 82=78; 83=79;
  84=80; 85=81
}

$color = @{
  43="0xff8f4b"; 46="0xff1741"
}

foreach ($r in $order.Keys) {
    $ws.Cells.Item($r, 28).Value = $order[$r]
    if ($color.ContainsKey($r)) {
        $ws.Cells.Item($r, 29).Value = $color[$r]
    } else {
        $ws.Cells.Item($r, 29).Value = "0x000000"
    }
}

# --- 6. Column widths for the two new columns -------------------------------
$w = $ws.Columns.Item(27).ColumnWidth()
$ws.Columns.Item(28).ColumnWidth = $w
$ws.Columns.Item(29).ColumnWidth = $w

# --- 7. View state: pane / selection matching the authored edit ------------
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("V3").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("AB33:AC33").Select()
